$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing row 2 (Ambiente/URL moved from "preproduccion" to "i-preproduccion", Plan bumped) ---

# Remove existing hyperlink(s) on sheet1 so we can cleanly re-create them
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws1.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws1.Range("L2").Value = "81 ver. 3"
$ws1.Range("N2").Value = 30990130

# --- Add new row 3, a copy of the (updated) row 2 but with a different Plan/NumeroDocumento ---

$ws1.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws1.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws1.Range("C3").Value = "su"
$ws1.Range("D3").Value = "silverarrow"
$ws1.Range("E3").Value = 3199801311
$ws1.Range("F3").Value = 2344
$ws1.Range("G3").Value = "Answer"
$ws1.Range("H3").Value = "Accidentes Personales"
$ws1.Range("I3").Value = "Cupón"
$ws1.Range("J3").Value = "No"
$ws1.Range("L3").Value = "Exe Autonomia ver. 2"
$ws1.Range("M3").Value = "D.N.I."
$ws1.Range("N3").Value = 30990131

# Match the number-formatting/font color already used on E2 (style index 2 in the template)
$ws1.Range("E3").Font.Color = $ws1.Range("E2").Font.Color

# Re-create the hyperlinks on B2 and B3 pointing at the new URL
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do")
$ws1.Range("B2").Style = "Hipervínculo"
$ws1.Range("B3").Style = "Hipervínculo"

# Update the sheet view: drop the frozen "topLeftCell" scroll and select B7 instead of N3
$ws1.Activate() | Out-Null
$ws1.Range("B7").Select() | Out-Null
